# Update column F (dSF) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    3  = 1
    5  = -2
    6  = -1
    8  = 7
    9  = -1
    10 = 0
    11 = -1
    12 = -2
    13 = 1
    14 = 5
    15 = 1
    16 = 0
    17 = 4
    18 = 1
    19 = -1
    21 = -4
    26 = -1
    30 = 2
    31 = 0
    32 = -1
    33 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
